$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # ALC
$ws1.Range("H62").Value = 2906.3
$ws1.Range("I62").Value = 2469.125
$ws1.Range("J62").Value = 4655
$ws1.Range("K62").Value = 2469.125
$ws1.Range("L62").Value = 4655
$ws1.Range("M62").Value = -1845.125
$ws1.Range("N62").Value = -5903

$ws1.Range("H65").Value = 2906.3
$ws1.Range("I65").Value = 2469.125
$ws1.Range("J65").Value = 4655
$ws1.Range("K65").Value = 12345.625
$ws1.Range("L65").Value = 23275
$ws1.Range("M65").Value = -9225.625
$ws1.Range("N65").Value = -29515

$ws1.Range("H70").Value = 4233.2144
$ws1.Range("I70").Value = 1000
$ws1.Range("K70").Value = 3000
$ws1.Range("M70").Value = -2730

$ws1.Range("H73").Value = 4233.2144
$ws1.Range("I73").Value = 1000
$ws1.Range("K73").Value = 3000
$ws1.Range("M73").Value = -2064

$ws1.Range("H132").Value = 22700.85
$ws1.Range("I132").Value = 28524
$ws1.Range("J132").Value = 1647.9231
$ws1.Range("K132").Value = 85572
$ws1.Range("L132").Value = 4943.7693
$ws1.Range("M132").Value = -83042
$ws1.Range("N132").Value = -10003.7693

$ws1.Range("H137").Value = 1339.4231
$ws1.Range("I137").Value = 1195.0513
$ws1.Range("J137").Value = 1772.5385
$ws1.Range("K137").Value = 3585.1539
$ws1.Range("L137").Value = 5317.6155
$ws1.Range("M137").Value = -1035.1539
$ws1.Range("N137").Value = -10417.6155

$ws2 = $wb.Worksheets.Item(2)  # ARM
$ws2.Range("H2").Value = 1028.6842
$ws2.Range("I2").Value = 1084.0333
$ws2.Range("J2").Value = 821.125
$ws2.Range("K2").Value = 1084.0333
$ws2.Range("L2").Value = 821.125
$ws2.Range("M2").Value = -971.0333000000001
$ws2.Range("N2").Value = -1047.125

$ws2.Range("H61").Value = 1192.976
$ws2.Range("I61").Value = 999.45715
$ws2.Range("J61").Value = 2235
$ws2.Range("K61").Value = 999.45715
$ws2.Range("L61").Value = 2235
$ws2.Range("M61").Value = -787.45715
$ws2.Range("N61").Value = -2659

$ws2.Range("H116").Value = 1028.6842
$ws2.Range("I116").Value = 1084.0333
$ws2.Range("J116").Value = 821.125
$ws2.Range("K116").Value = 1084.0333
$ws2.Range("L116").Value = 821.125
$ws2.Range("M116").Value = 1209.9667
$ws2.Range("N116").Value = -5409.125

$ws2.Range("H119").Value = 45000
$ws2.Range("J119").Value = 45000
$ws2.Range("L119").Value = 45000
$ws2.Range("N119").Value = -54676

$ws2.Range("H132").Value = 1251.35
$ws2.Range("I132").Value = 935.08887
$ws2.Range("K132").Value = 2805.26661
$ws2.Range("M132").Value = -275.2666100000001

$ws2.Range("H136").Value = 1192.976
$ws2.Range("I136").Value = 999.45715
$ws2.Range("J136").Value = 2235
$ws2.Range("K136").Value = 2998.37145
$ws2.Range("L136").Value = 6705
$ws2.Range("M136").Value = -448.3714499999996
$ws2.Range("N136").Value = -11805

$ws3 = $wb.Worksheets.Item(3)  # BSM
$ws3.Range("H3").Value = 1028.6842
$ws3.Range("I3").Value = 1084.0333
$ws3.Range("J3").Value = 821.125
$ws3.Range("K3").Value = 1084.0333
$ws3.Range("L3").Value = 821.125
$ws3.Range("M3").Value = -970.0333000000001
$ws3.Range("N3").Value = -1049.125

$ws4 = $wb.Worksheets.Item(4)  # CRP
$ws4.Range("H31").Value = 2272.257
$ws4.Range("I31").Value = 1196.4615
$ws4.Range("J31").Value = 3625.6775
$ws4.Range("K31").Value = 1196.4615
$ws4.Range("L31").Value = 3625.6775
$ws4.Range("M31").Value = -901.4614999999999
$ws4.Range("N31").Value = -4215.6775

$ws4.Range("H34").Value = 2272.257
$ws4.Range("I34").Value = 1196.4615
$ws4.Range("J34").Value = 3625.6775
$ws4.Range("K34").Value = 1196.4615
$ws4.Range("L34").Value = 3625.6775
$ws4.Range("M34").Value = -994.4614999999999
$ws4.Range("N34").Value = -4029.6775

$ws4.Range("H132").Value = 1689.0204
$ws4.Range("I132").Value = 1270.069
$ws4.Range("J132").Value = 2296.5
$ws4.Range("K132").Value = 3810.207
$ws4.Range("L132").Value = 6889.5
$ws4.Range("M132").Value = -1280.207
$ws4.Range("N132").Value = -11949.5

$ws4.Range("H134").Value = 1328.2877
$ws4.Range("I134").Value = 1236
$ws4.Range("J134").Value = 1717.2142
$ws4.Range("K134").Value = 3708
$ws4.Range("L134").Value = 5151.642599999999
$ws4.Range("M134").Value = -1173
$ws4.Range("N134").Value = -10221.6426

$ws5 = $wb.Worksheets.Item(5)  # CUL
$ws5.Range("H122").Value = 1753.72
$ws5.Range("I122").Value = 471
$ws5.Range("J122").Value = 3386.2727
$ws5.Range("K122").Value = 4239
$ws5.Range("L122").Value = 30476.4543
$ws5.Range("M122").Value = -1789
$ws5.Range("N122").Value = -35376.4543

$ws5.Range("H131").Value = 2881.2876
$ws5.Range("I131").Value = 533.5454999999999
$ws5.Range("J131").Value = 3297.8225
$ws5.Range("K131").Value = 1600.6365
$ws5.Range("L131").Value = 9893.467500000001
$ws5.Range("M131").Value = 3439.3635
$ws5.Range("N131").Value = -19973.4675

$ws5.Range("H137").Value = 2184.95
$ws5.Range("I137").Value = 2411
$ws5.Range("J137").Value = 2000
$ws5.Range("K137").Value = 7233
$ws5.Range("L137").Value = 6000
$ws5.Range("M137").Value = -2133
$ws5.Range("N137").Value = -16200

$ws6 = $wb.Worksheets.Item(6)  # GSM
$ws6.Range("H132").Value = 2210.879
$ws6.Range("I132").Value = 1833.0869
$ws6.Range("J132").Value = 3079.8
$ws6.Range("K132").Value = 5499.2607
$ws6.Range("L132").Value = 9239.400000000001
$ws6.Range("M132").Value = -2969.2607
$ws6.Range("N132").Value = -14299.4

$ws7 = $wb.Worksheets.Item(7)  # LTW
$ws7.Range("H132").Value = 2071.78
$ws7.Range("I132").Value = 1905.3733
$ws7.Range("J132").Value = 2571
$ws7.Range("K132").Value = 5716.1199
$ws7.Range("L132").Value = 7713
$ws7.Range("M132").Value = -3186.1199
$ws7.Range("N132").Value = -12773

$ws8 = $wb.Worksheets.Item(8)  # WVR
$ws8.Range("H119").Value = 27200
$ws8.Range("J119").Value = 27200
$ws8.Range("L119").Value = 27200
$ws8.Range("N119").Value = -36876

$ws8.Range("H132").Value = 1629.4509
$ws8.Range("I132").Value = 1467.9487
$ws8.Range("J132").Value = 2154.3333
$ws8.Range("K132").Value = 4403.8461
$ws8.Range("L132").Value = 6462.999899999999
$ws8.Range("M132").Value = -1873.8461
$ws8.Range("N132").Value = -11522.9999

$ws8.Range("H136").Value = 5209409
$ws8.Range("I136").Value = 8333841.5
$ws8.Range("J136").Value = 2021.0416
$ws8.Range("K136").Value = 25001524.5
$ws8.Range("L136").Value = 6063.1248
$ws8.Range("M136").Value = -24998974.5
$ws8.Range("N136").Value = -11163.1248
